$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row before row 259 by duplicating the current row 259
# (same Mercado/Region/Producto/Variedad template) and shifting the rest of
# the table down by one row.
$ws.Rows.Item(259).Copy()
$ws.Rows.Item(259).Insert()

# Overwrite the new row 259 with the new record's values.
$ws.Range("D259").Value2 = 44875
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value2 = 400
$ws.Range("N259").Value2 = 17000
$ws.Range("O259").Value2 = 18000
$ws.Range("P259").Value2 = 17500
$ws.Range("S259").Value2 = 1167
